# Auto-generated edit script applying the Tiamat_Profits.xlsx diff
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N) across all 8 job sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 359.52
$ws.Range("I15").Value = 359.52
$ws.Range("K15").Value = 1078.56
$ws.Range("M15").Value = -909.5599999999999
# Row 98
$ws.Range("H98").Value = 7937616.5
$ws.Range("J98").Value = 1597.6666
$ws.Range("L98").Value = 1597.6666
$ws.Range("N98").Value = -4593.6666
# Row 103
$ws.Range("H103").Value = 772377.3
$ws.Range("I103").Value = 670
$ws.Range("J103").Value = 1323596.9
$ws.Range("K103").Value = 2010
$ws.Range("L103").Value = 3970790.7
$ws.Range("M103").Value = -1424
$ws.Range("N103").Value = -3971962.7
# Row 106
$ws.Range("H106").Value = 46876.25
$ws.Range("I106").Value = 70402
$ws.Range("J106").Value = 7666.6665
$ws.Range("K106").Value = 70402
$ws.Range("L106").Value = 7666.6665
$ws.Range("M106").Value = -69771
$ws.Range("N106").Value = -8928.666499999999
# Row 122
$ws.Range("H122").Value = 7937616.5
$ws.Range("J122").Value = 1597.6666
$ws.Range("L122").Value = 4792.9998
$ws.Range("N122").Value = -9692.9998
# Row 124
$ws.Range("H124").Value = 34926.668
$ws.Range("J124").Value = 34926.668
$ws.Range("L124").Value = 34926.668
$ws.Range("N124").Value = -44746.668
# Row 132
$ws.Range("H132").Value = 172094.56
$ws.Range("I132").Value = 2488.389
$ws.Range("K132").Value = 7465.167
$ws.Range("M132").Value = -4935.167

$ws = $wb.Worksheets.Item("ARM")
# Row 10
$ws.Range("H10").Value = 75004
$ws.Range("I10").Value = 5000
$ws.Range("J10").Value = 92505
$ws.Range("K10").Value = 5000
$ws.Range("L10").Value = 92505
$ws.Range("M10").Value = -4830
$ws.Range("N10").Value = -92845
# Row 32
$ws.Range("H32").Value = 3252.13
$ws.Range("I32").Value = 2994.819
$ws.Range("J32").Value = 7283.3335
$ws.Range("K32").Value = 2994.819
$ws.Range("L32").Value = 7283.3335
$ws.Range("M32").Value = -2707.819
$ws.Range("N32").Value = -7857.3335
# Row 97
$ws.Range("H97").Value = 880.2083
$ws.Range("I97").Value = 460.68182
$ws.Range("J97").Value = 5495
$ws.Range("K97").Value = 460.68182
$ws.Range("L97").Value = 5495
$ws.Range("M97").Value = 35.31817999999998
$ws.Range("N97").Value = -6487

$ws = $wb.Worksheets.Item("BSM")
# Row 16
$ws.Range("H16").Value = 34503
$ws.Range("I16").Value = 10336.333
$ws.Range("J16").Value = 58669.668
$ws.Range("K16").Value = 10336.333
$ws.Range("L16").Value = 58669.668
$ws.Range("M16").Value = -10166.333
$ws.Range("N16").Value = -59009.668
# Row 86
$ws.Range("H86").Value = 413459.2
$ws.Range("I86").Value = 1671.3334
$ws.Range("J86").Value = 1401750
$ws.Range("K86").Value = 1671.3334
$ws.Range("L86").Value = 1401750
$ws.Range("M86").Value = -548.3334
$ws.Range("N86").Value = -1403996
# Row 89
$ws.Range("H89").Value = 413459.2
$ws.Range("I89").Value = 1671.3334
$ws.Range("J89").Value = 1401750
$ws.Range("K89").Value = 8356.666999999999
$ws.Range("L89").Value = 7008750
$ws.Range("M89").Value = -2740.666999999999
$ws.Range("N89").Value = -7019982
# Row 94
$ws.Range("H94").Value = 3709.9048
$ws.Range("I94").Value = 693.26666
$ws.Range("J94").Value = 11251.5
$ws.Range("K94").Value = 693.26666
$ws.Range("L94").Value = 11251.5
$ws.Range("M94").Value = -242.26666
$ws.Range("N94").Value = -12153.5
# Row 99
$ws.Range("H99").Value = 1378.6364
$ws.Range("I99").Value = 1172.3529
$ws.Range("K99").Value = 1172.3529
$ws.Range("M99").Value = 325.6470999999999
# Row 105
$ws.Range("H105").Value = 1593962
$ws.Range("I105").Value = 3981005
$ws.Range("J105").Value = 2600
$ws.Range("K105").Value = 3981005
$ws.Range("L105").Value = 2600
$ws.Range("M105").Value = -3979258
$ws.Range("N105").Value = -6094

$ws = $wb.Worksheets.Item("CRP")
# Row 50
$ws.Range("H50").Value = 9884.546
$ws.Range("J50").Value = 9884.546
$ws.Range("L50").Value = 9884.546
$ws.Range("N50").Value = -11134.546
# Row 141
$ws.Range("H141").Value = 27305.555
$ws.Range("I141").Value = 6545.3335
$ws.Range("J141").Value = 33237.047
$ws.Range("K141").Value = 6545.3335
$ws.Range("L141").Value = 33237.047
$ws.Range("M141").Value = -1365.3335
$ws.Range("N141").Value = -43597.047

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 1236.5714
$ws.Range("I5").Value = 845.8889
$ws.Range("J5").Value = 1650.2354
$ws.Range("K5").Value = 2537.6667
$ws.Range("L5").Value = 4950.706200000001
$ws.Range("M5").Value = -2425.6667
$ws.Range("N5").Value = -5174.706200000001
# Row 61
$ws.Range("H61").Value = 110
$ws.Range("I61").Value = 100
$ws.Range("J61").Value = 150
$ws.Range("K61").Value = 300
$ws.Range("L61").Value = 450
$ws.Range("M61").Value = -85
$ws.Range("N61").Value = -880
# Row 113
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").ClearContents()
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = 0
# Row 131
$ws.Range("H131").Value = 294884.72
$ws.Range("J131").Value = 323379.7
$ws.Range("L131").Value = 970139.1000000001
$ws.Range("N131").Value = -980219.1000000001
# Row 135
$ws.Range("H135").Value = 1236.5714
$ws.Range("I135").Value = 845.8889
$ws.Range("J135").Value = 1650.2354
$ws.Range("K135").Value = 7613.0001
$ws.Range("L135").Value = 14852.1186
$ws.Range("M135").Value = -5078.0001
$ws.Range("N135").Value = -19922.1186

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 8417.883
$ws.Range("I80").Value = 2911.5557
$ws.Range("J80").Value = 14612.5
$ws.Range("K80").Value = 2911.5557
$ws.Range("L80").Value = 14612.5
$ws.Range("M80").Value = -1913.5557
$ws.Range("N80").Value = -16608.5
# Row 83
$ws.Range("H83").Value = 8417.883
$ws.Range("I83").Value = 2911.5557
$ws.Range("J83").Value = 14612.5
$ws.Range("K83").Value = 14557.7785
$ws.Range("L83").Value = 73062.5
$ws.Range("M83").Value = -9565.7785
$ws.Range("N83").Value = -83046.5
# Row 97
$ws.Range("H97").Value = 1096.625
$ws.Range("I97").Value = 1075
$ws.Range("J97").Value = 1118.25
$ws.Range("K97").Value = 1075
$ws.Range("L97").Value = 1118.25
$ws.Range("M97").Value = -579
$ws.Range("N97").Value = -2110.25
# Row 122
$ws.Range("H122").Value = 2871.8462
$ws.Range("I122").Value = 2515.625
$ws.Range("J122").Value = 3441.8
$ws.Range("K122").Value = 7546.875
$ws.Range("L122").Value = 10325.4
$ws.Range("M122").Value = -5096.875
$ws.Range("N122").Value = -15225.4

$ws = $wb.Worksheets.Item("LTW")
# Row 82
$ws.Range("H82").Value = 1522.0834
$ws.Range("I82").Value = 1071.2727
$ws.Range("J82").Value = 2230.5
$ws.Range("K82").Value = 1071.2727
$ws.Range("L82").Value = 2230.5
$ws.Range("M82").Value = -710.2727
$ws.Range("N82").Value = -2952.5
# Row 85
$ws.Range("H85").Value = 1522.0834
$ws.Range("I85").Value = 1071.2727
$ws.Range("J85").Value = 2230.5
$ws.Range("K85").Value = 1071.2727
$ws.Range("L85").Value = 2230.5
$ws.Range("M85").Value = 176.7273
$ws.Range("N85").Value = -4726.5

$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 1436.4482
$ws.Range("I81").Value = 1549.2778
$ws.Range("J81").Value = 1251.8182
$ws.Range("K81").Value = 3098.5556
$ws.Range("L81").Value = 2503.6364
$ws.Range("M81").Value = -2037.5556
$ws.Range("N81").Value = -4625.636399999999
# Row 84
$ws.Range("H84").Value = 1436.4482
$ws.Range("I84").Value = 1549.2778
$ws.Range("J84").Value = 1251.8182
$ws.Range("K84").Value = 15492.778
$ws.Range("L84").Value = 12518.182
$ws.Range("M84").Value = -10188.778
$ws.Range("N84").Value = -23126.182
# Row 125
$ws.Range("H125").Value = 24390.445
$ws.Range("J125").Value = 24390.445
$ws.Range("L125").Value = 24390.445
$ws.Range("N125").Value = -34230.445
# Row 132
$ws.Range("H132").Value = 4172.121
$ws.Range("I132").Value = 779.3913
$ws.Range("J132").Value = 11975.4
$ws.Range("K132").Value = 2338.1739
$ws.Range("L132").Value = 35926.2
$ws.Range("M132").Value = 191.8261000000002
$ws.Range("N132").Value = -40986.2
